$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Student 1"
$ws.Range("B2").Value = "student1@example.com"
$ws.Range("C2").Value = "p1,p2,p3"

# "123" looks numeric, so a plain .Value assignment would store it as a
# number (no shared-string entry, t="s" would be lost). Build it as a
# text-formula result in a scratch cell, then paste only the value back
# into D2 so it lands as a genuine text/shared-string cell without
# picking up any NumberFormat-driven style changes.
$ws.Range("Z1").Formula = "=""123"""
$ws.Range("Z1").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("E2").Value = "Inactive"

$ws.Range("F11").Select()
